# Add functions checkAndAddBlankShortcuts and rearrangeGroupIds
# -> fills in the groupBg color for the "SUPP" group (rows 12-17) and
#    moves the stray shortcutBg value in row 16 into its correct groupBg column,
#    then updates that color's hex value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 currently has its color value mis-placed in column F (shortcutBg)
# instead of column C (groupBg). Move it over.
$strayValue = $ws.Range("F16").Value2
$ws.Range("C16").Value2 = $strayValue
$ws.Range("F16").ClearContents()

# Backfill the groupBg color for the rest of the SUPP group (group-2) rows,
# reusing the same shared string as the moved value above.
$ws.Range("C12").Value2 = $strayValue
$ws.Range("C13").Value2 = $strayValue
$ws.Range("C14").Value2 = $strayValue
$ws.Range("C15").Value2 = $strayValue
$ws.Range("C17").Value2 = $strayValue

# Update the color itself everywhere it is used.
$ws.Cells.Replace("#7D3C98", "#e8a5a0")

# Restore the active selection to match where the edit left off.
$ws.Range("M17").Select()
